$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.321.67"
$ws.Range("E2").Value = "  -3.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.007.98"
$ws.Range("E3").Value = "  -4.93%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.49"
$ws.Range("E5").Value = "  -2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.86"
$ws.Range("E6").Value = "  -7.14%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.010.32"
$ws.Range("E9").Value = "  -4.78%  "

$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("E11").Value = "  -6.18%  "

$ws.Range("E12").Value = "  -3.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.534.49"
$ws.Range("E13").Value = "  -4.83%  "

$ws.Range("E14").Value = "  -3.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.466.70"
$ws.Range("E15").Value = "  -2.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.84"
$ws.Range("E16").Value = "  -5.48%  "

$ws.Range("E17").Value = "  -3.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.011.76"
$ws.Range("E18").Value = "  -4.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "394.45"
$ws.Range("E19").Value = "  -4.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.13"
$ws.Range("E20").Value = "  -2.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.03"
$ws.Range("E21").Value = "  -4.43%  "

$ws.Range("E22").Value = "  -6.09%  "

$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.29"
$ws.Range("E24").Value = "  -3.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.470"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("E26").Value = "  -7.11%  "

$ws.Range("E27").Value = "  -5.51%  "

$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("E29").Value = "  -3.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -3.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.60"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.01"
$ws.Range("E33").Value = "  +4.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  -2.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.05"
$ws.Range("E35").Value = "  -3.81%  "

$ws.Range("E36").Value = "  -3.18%  "

$ws.Range("E37").Value = "  -2.95%  "

$ws.Range("E38").Value = "  -6.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.460.60"
$ws.Range("E39").Value = "  -10.00%  "

$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.63"
$ws.Range("E41").Value = "  -3.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.62"
$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("E43").Value = "  -5.42%  "

$ws.Range("E44").Value = "  -4.44%  "

$ws.Range("E45").Value = "  -4.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("E47").Value = "  -9.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.92"
$ws.Range("E48").Value = "  -5.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0957"
$ws.Range("E49").Value = "  -2.42%  "

$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "265.36"
$ws.Range("E51").Value = "  -7.37%  "
